$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $cellRef, $value)
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-CellText $ws 'D2' '245.12'
Set-CellText $ws 'E2' '-0.49%'
Set-CellText $ws 'E3' '-4.51%'
Set-CellText $ws 'D4' '5.234'
Set-CellText $ws 'D5' '0.05694'
Set-CellText $ws 'E5' '-0.57%'
Set-CellText $ws 'D6' '6.616'
Set-CellText $ws 'E6' '0.34%'
Set-CellText $ws 'D7' '3.196'
Set-CellText $ws 'E7' '3.17%'
Set-CellText $ws 'D8' '0.8501'
Set-CellText $ws 'E8' '-0.69%'
Set-CellText $ws 'D9' '0.8576'
Set-CellText $ws 'E9' '-1.27%'
Set-CellText $ws 'B10' 'One'
Set-CellText $ws 'C10' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-CellText $ws 'D10' '0.01001'
Set-CellText $ws 'E10' '1,563.63%'
Set-CellText $ws 'B11' 'WazirX'
Set-CellText $ws 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-CellText $ws 'D11' '0.1370'
Set-CellText $ws 'E11' '0.31%'
Set-CellText $ws 'B12' 'MandalaExchangeToken'
Set-CellText $ws 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-CellText $ws 'D12' '0.07023'
Set-CellText $ws 'E12' '-0.65%'
Set-CellText $ws 'B13' 'BitrueCoin'
Set-CellText $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-CellText $ws 'D13' '0.03136'
Set-CellText $ws 'E13' '6.98%'
Set-CellText $ws 'B14' 'BitMartToken'
Set-CellText $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-CellText $ws 'D14' '0.09210'
Set-CellText $ws 'E14' '-1.86%'
Set-CellText $ws 'B15' 'BitForexToken'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-CellText $ws 'D15' '0.001542'
Set-CellText $ws 'E15' '1.04%'
Set-CellText $ws 'D16' '0.005896'
Set-CellText $ws 'E16' '-4.99%'
Set-CellText $ws 'E17' '-0.02%'
Set-CellText $ws 'E19' '0.42%'
Set-CellText $ws 'D20' '0.03259'
Set-CellText $ws 'E20' '-5.33%'
Set-CellText $ws 'D21' '0.1288'
Set-CellText $ws 'E21' '-2.22%'
Set-CellText $ws 'D22' '3.492'
Set-CellText $ws 'E22' '0.68%'
Set-CellText $ws 'E23' '-2.20%'
Set-CellText $ws 'D24' '0.1379'
Set-CellText $ws 'E24' '-0.06%'
Set-CellText $ws 'D25' '0.001217'
Set-CellText $ws 'E25' '-1.26%'
Set-CellText $ws 'D26' '0.004140'
Set-CellText $ws 'E26' '-17.49%'
Set-CellText $ws 'D27' '0.0001199'
Set-CellText $ws 'E27' '-0.88%'
Set-CellText $ws 'D28' '0.0001448'
Set-CellText $ws 'D40' '0.03771'
Set-CellText $ws 'E40' '0.53%'
Set-CellText $ws 'D41' '0.1063'
Set-CellText $ws 'E41' '-0.79%'
Set-CellText $ws 'D42' '0.003736'
Set-CellText $ws 'E42' '-35.24%'
Set-CellText $ws 'D43' '0.002298'
Set-CellText $ws 'E43' '14.94%'
Set-CellText $ws 'D44' '0.009152'
Set-CellText $ws 'E44' '-4.48%'
Set-CellText $ws 'D45' '0.00005285'
Set-CellText $ws 'E45' '1.17%'
Set-CellText $ws 'E46' '-0.05%'
Set-CellText $ws 'D47' '0.1049'
Set-CellText $ws 'E47' '62.20%'
Set-CellText $ws 'D48' '0.002437'
Set-CellText $ws 'E48' '-3.33%'
Set-CellText $ws 'D49' '0.00002099'
Set-CellText $ws 'E49' '-0.05%'
Set-CellText $ws 'D50' '0.0001999'
Set-CellText $ws 'E50' '-0.05%'
